# Update Data by bot, scripted by HH
# Row 5 holds the 2020 Q3 report figures for security 003028; restate it
# to the corresponding 2019 Q3 figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# REPORTDATE changed from 2020-09-30 to 2019-09-30
$ws.Range("H5").Value = "2019-09-30 00:00:00"

$ws.Range("I5").Value = 0.97
$ws.Range("J5").ClearContents()

$ws.Range("K5").Value = 509598813.87
$ws.Range("L5").Value = 79455072.47

$ws.Range("M5").ClearContents()
$ws.Range("N5").ClearContents()
$ws.Range("O5").ClearContents()
$ws.Range("P5").ClearContents()
$ws.Range("Q5").ClearContents()

$ws.Range("R5").Value = 27.0353207426

# AB5 and AE5 look like numbers ("0" and "2019") so force them to be
# stored as text, matching the original inlineStr type, then reset the
# cell style back to Normal so no stray number-format style is left on it.
$ws.Range("AB5").NumberFormat = "@"
$ws.Range("AB5").Value = "0"
$ws.Range("AB5").Style = "Normal"

$ws.Range("AC5").Value = "2019Q3"
$ws.Range("AD5").Value = "2019年 三季报"

$ws.Range("AE5").NumberFormat = "@"
$ws.Range("AE5").Value = "2019"
$ws.Range("AE5").Style = "Normal"
